$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update logged time entries (Start/End/Break times) ---

# Row 6 (2020-10-07ish): End time (C6) pushed later
$ws.Range("C6").Value = 0.65806712962962965

# Row 19: End time (C19) pushed later
$ws.Range("C19").Value = 0.7901273148148148

# Row 29: Break duration (D29) increased from 20 to 30 minutes
$ws.Range("D29").Value = 0.020833333333333332

# Row 30: previously an empty/untouched day - now has logged Start/End time and an activity note
$ws.Range("B30").Value = 0.38245370370370368
$ws.Range("C30").Value = 0.79780092592592589

$newActivity = "* Added runtime tool implementation`n* Added custom inspector for runtime MonoBehaviour`n* Added small demo project showcasing the runtime functionality of the tool`n* Added alerts for graph file being deleted`n* Added functionality to attempt to recover graph file if it ended up as null`n* Fixed bug where runtime 'Next' property of Nodes wasn't set properly`n* Fixed small styling bug"
$ws.Range("F30").Value = $newActivity
$ws.Range("F30").WrapText = $true
$ws.Rows("30").RowHeight = 15

# --- Widen the Activity column (F) to better fit the new, longer log entry ---
$ws.Columns("F").ColumnWidth = 36.7

# --- Move the active selection to C31 (next empty entry to fill in) ---
$ws.Range("C31").Select() | Out-Null
